# Auto-update Cloudflare export
# Replace the legacy "{'flatten_cname': False}" settings value with "{}"
# in column H for the DNS records that previously had flatten_cname
# explicitly set to False.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$rows = @(3,13,14,15,16,26,31,37,38,39,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,83,84,85,86,87,88,90,91,92,93,94,110,111,112,113,114,115,116,127,128,141)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 8)  # Column H = 8
    if ($cell.Value2 -eq "{'flatten_cname': False}") {
        $cell.Value = "{}"
    }
}
